$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cryptocurrency Price/Volume columns (D and E) store plain text,
# including numeric-looking values (e.g. "67.72"). Force each target
# cell to Text format first so Excel does not auto-convert the new
# value to a number, matching the original inline-string text cells.
# (NumberFormat/Style must be set per-cell, not via a multi-area
# Range, since multi-area assignment only reliably affects the
# first area in this engine.)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.222.03"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.648.19"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "218.28"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "20.23"
$ws.Range("E10").Value = "  +3.17%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D13").Value = "1.639.49"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "0.538"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "67.72"
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").Value = "27.190.57"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "220.60"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "6.76"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").Value = "4.44"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "148.62"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "1.274.01"
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").Value = "0.845"
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "2.27"
$ws.Range("E41").Value = "  +7.90%  "
$ws.Range("D42").Value = "0.810"
$ws.Range("E42").Value = "  +0.28%  "
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "1.788.16"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "62.86"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").Value = "92.36"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("E48").Value = "  +17.20%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  -0.13%  "

# Restore the default (unstyled) cell style now that the text values
# are safely stored, so no stray formatting is introduced.
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Style = "Normal"
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Style = "Normal"
